$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "All input"

# Add the new "1 input" sheet right after "All input"
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "1 input"

# ---------------------------------------------------------------------------
# Update "All input" sheet data
# ---------------------------------------------------------------------------

# Row 6: no new data, but fix up the odd bold formatting on J6/K6 so it matches
# the rest of the Train/Test RMSE columns (J3:J5 / K3:K5).
$ws1.Range("J3").Copy()
$ws1.Range("J6").PasteSpecial(-4122)
$ws1.Range("K3").Copy()
$ws1.Range("K6").PasteSpecial(-4122)

# Row 7: change model name A7 from LSTM to GRU, fill J7:M7
$ws1.Range("A7").Value = "GRU - scale min max"
$ws1.Range("J7").Value = 0.070327
$ws1.Range("K7").Value = 0.069386
$ws1.Range("L7").Value = 0.055053
$ws1.Range("M7").Value = 0.040517

# Row 8: fill entire row with new data (was blank)
$ws1.Range("A8").Value = "GRU - scale min max"
$ws1.Range("B8").Value = 0.001
$ws1.Range("C8").Value = 20
$ws1.Range("D8").Value = 4
$ws1.Range("E8").Value = 128
$ws1.Range("F8").Value = 0.8
$ws1.Range("G8").Value = 4
$ws1.Range("H8").Value = 32
$ws1.Range("I8").Value = 18
$ws1.Range("J8").Value = 0.069891
$ws1.Range("K8").Value = 0.066648
$ws1.Range("L8").Value = 0.055179
$ws1.Range("M8").Value = 0.039709

# Row 9: brand new row
$ws1.Range("A9").Value = "GRU - scale min max"
$ws1.Range("B9").Value = 0.001
$ws1.Range("C9").Value = 20
$ws1.Range("D9").Value = 4
$ws1.Range("E9").Value = 256
$ws1.Range("F9").Value = 0.8
$ws1.Range("G9").Value = 4
$ws1.Range("H9").Value = 32
$ws1.Range("I9").Value = 18
$ws1.Range("J9").Value = 0.063453
$ws1.Range("K9").Value = 0.064912
$ws1.Range("L9").Value = 0.053531
$ws1.Range("M9").Value = 0.041204

# Make sure A9:I9 take the same row-level formatting as the rest of the table
# (rather than falling back to the plain column default format).
$ws1.Range("A3:I3").Copy()
$ws1.Range("A9:I9").PasteSpecial(-4122)

# J7:M7 and J8:M8 should use the same formatting as the other Train/Test RMSE
# and 1-step-ahead RMSE cells.
$ws1.Range("J3").Copy()
$ws1.Range("J7").PasteSpecial(-4122)
$ws1.Range("J8").PasteSpecial(-4122)
$ws1.Range("K3").Copy()
$ws1.Range("K7").PasteSpecial(-4122)
$ws1.Range("K8").PasteSpecial(-4122)
$ws1.Range("L5").Copy()
$ws1.Range("L3").PasteSpecial(-4122)
$ws1.Range("L4").PasteSpecial(-4122)
$ws1.Range("L7").PasteSpecial(-4122)
$ws1.Range("L8").PasteSpecial(-4122)
$ws1.Range("M5").Copy()
$ws1.Range("M3").PasteSpecial(-4122)
$ws1.Range("M4").PasteSpecial(-4122)
$ws1.Range("M7").PasteSpecial(-4122)
$ws1.Range("M8").PasteSpecial(-4122)

# Row heights: 6,7,8 become 27 (matching rows 4 & 5)
$ws1.Rows.Item(6).RowHeight = 27
$ws1.Rows.Item(7).RowHeight = 27
$ws1.Rows.Item(8).RowHeight = 27

# Selection / active cell bookkeeping
$ws1.Range("B11").Select()
$ws2.Activate()

Write-Host "Done"
